$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 new rows before the old "Total Hours" row (row 11) -------
$ws.Range("A11:A14").EntireRow.Insert()

# Copy the formatting (borders/fill/font/alignment) from row 10's data
# cells down into the freshly inserted rows 11-14 so they keep the same
# "activity line" look (style indices 8 / 9) instead of Excel's blank
# unformatted default.
$ws.Range("B10:C10").Copy()
$ws.Range("B11:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Update the header period cell --------------------------------------
$ws.Range("A5").Value = "تیر 99"

# --- 3. Update the activity / hours / task rows -----------------------------
$ws.Range("B6").Value = "* Documents (Read/Write)"
$ws.Range("C6").Value = 1
$ws.Range("E6").Value = "• Render on GPU"

$ws.Range("B7").Value = "* Segmentation (Lung/Airway)"
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = "• Segmentation (Lung/Airway)       (not integrated)"

$ws.Range("B8").Value = "* CenterLine Extraction"
$ws.Range("C8").Value = 9
$ws.Range("E8").Value = "• Centerline extraction                    (not integrated)"

$ws.Range("B9").Value = "* Registration"
$ws.Range("C9").Value = 6
# E9 was empty (no style) before - pull E8's formatting (thin border box,
# the same one E6/E7/E8 use) across before writing its text.
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E9").Value = "• Registration (CPD)                         (not integrated)"

$ws.Range("B10").Value = "* GPU Support"
$ws.Range("C10").Value = 2

$ws.Range("B11").Value = "* GUI"
$ws.Range("C11").Value = 4

$ws.Range("B12").Value = "* Support for VTK files"
$ws.Range("C12").Value = 3

$ws.Range("B13").Value = "* Integration"
$ws.Range("C13").Value = 3

$ws.Range("B14").Value = "* Meetings"
$ws.Range("C14").Value = 2

# --- 4. Totals block (now shifted down to rows 15-17) -----------------------
# B15/C16/D16/C17 already hold the correct text/value after the row-insert
# shifted them down from B11/C12/D12/C13 - leave them untouched so their
# original cell styles (incl. the quotePrefix "@..." style) survive, and
# only touch the two cells whose formula actually changes.
$ws.Range("C15").Formula = "=SUM(C6:C14)"
$ws.Range("D17").Formula = "=C15"

# --- 5. Column E width & selection ------------------------------------------
# ColumnWidth (character units) and the stored OOXML <col width> differ by
# the sheet's max-digit-width padding (~0.8333 chars for Calibri 11) - back
# that off so the saved width lands on exactly 42.
$ws.Columns.Item(5).ColumnWidth = 41.1666666666667
$ws.Range("E13").Select()

$wb.Application.CalculateFull()
